$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = Get-Date -Year 2021 -Month 10 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Range("J2").Value = 150

# Row 3
$ws.Range("D3").Value = Get-Date -Year 2021 -Month 10 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("J3").Value = 250

# Row 5
$ws.Range("D5").Value = Get-Date -Year 2021 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = 800
$ws.Range("N5").Value = '$/kilo (volumen en unidades)'
$ws.Range("O5").Value = 'Perú'
$ws.Range("P5").Value = 800

# Row 6
$ws.Range("D6").Value = Get-Date -Year 2021 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 800
$ws.Range("L6").Value = 800
$ws.Range("M6").Value = 800
$ws.Range("N6").Value = '$/kilo (volumen en unidades)'
$ws.Range("O6").Value = 'Perú'
$ws.Range("P6").Value = 800

# Row 7
$ws.Range("D7").Value = Get-Date -Year 2021 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2500
$ws.Range("O7").Value = 'Perú'
$ws.Range("P7").Value = 2500

# Row 8
$ws.Range("D8").Value = Get-Date -Year 2021 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 800
$ws.Range("N8").Value = '$/kilo (volumen en unidades)'
$ws.Range("O8").Value = 'Perú'
$ws.Range("P8").Value = 800

# Row 9
$ws.Range("D9").Value = Get-Date -Year 2020 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("I9").Value = 'Extra'
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 3500
$ws.Range("L9").Value = 3500
$ws.Range("M9").Value = 3500
$ws.Range("N9").Value = '$/unidad'
$ws.Range("O9").Value = 'Región de O''Higgins'
$ws.Range("P9").Value = 3500

# Row 10
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("P10").Value = 3000

# Row 11
$ws.Range("D11").Value = Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("I11").Value = 'Extra'
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2500
$ws.Range("P11").Value = 2500

# Row 12
$ws.Range("D12").Value = Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("J12").Value = 280
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 2000
$ws.Range("N12").Value = '$/unidad'
$ws.Range("O12").Value = 'Región de O''Higgins'
$ws.Range("P12").Value = 2000

# Row 13
$ws.Range("D13").Value = Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H13").Value = 'Americana O Klondike'
$ws.Range("I13").Value = 'Extra'
$ws.Range("J13").Value = 340
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2500
$ws.Range("N13").Value = '$/unidad'
$ws.Range("O13").Value = 'Región de O''Higgins'
$ws.Range("P13").Value = 2500

# Row 14
$ws.Range("D14").Value = Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H14").Value = 'Americana O Klondike'
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 2000
$ws.Range("O14").Value = 'Región de O''Higgins'
$ws.Range("P14").Value = 2000

# Row 15
$ws.Range("D15").Value = Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H15").Value = 'Americana O Klondike'
$ws.Range("I15").Value = 'Segunda'
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 1500
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = 1500
$ws.Range("N15").Value = '$/unidad'
$ws.Range("O15").Value = 'Región de O''Higgins'
$ws.Range("P15").Value = 1500

# Row 16
$ws.Range("D16").Value = Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H16").Value = 'Americana O Klondike'
$ws.Range("I16").Value = 'Tercera'
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 1000
$ws.Range("P16").Value = 1000

# Row 17
$ws.Range("D17").Value = Get-Date -Year 2021 -Month 10 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 800
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = 800
$ws.Range("N17").Value = '$/kilo (volumen en unidades)'
$ws.Range("O17").Value = 'Perú'
$ws.Range("P17").Value = 800

# Row 19
$ws.Range("D19").Value = Get-Date -Year 2021 -Month 4 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2500
$ws.Range("O19").Value = 'Perú'
$ws.Range("P19").Value = 2500

# Row 20
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = 5000
$ws.Range("P20").Value = 5000

# Row 21
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 560
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = 3000
$ws.Range("P21").Value = 3000

# Row 22
$ws.Range("D22").Value = Get-Date -Year 2020 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("I22").Value = 'Tercera'
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 2000
$ws.Range("O22").Value = 'Región de O''Higgins'
$ws.Range("P22").Value = 2000
